$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") for rows 2-20 is being updated from serial date
# 45224 (2023-10-25) to 45233 (2023-11-03).
$ws.Range("C2:C20").Value = 45233
